# Updates the multiplication problems in the table to the new master set.
$d = $word.ActiveDocument

$replacements = @(
    @("11×44=", "70×90="),
    @("81×59=", "32×29="),
    @("87×65=", "70×70="),
    @("66×20=", "56×41="),
    @("25×43=", "87×76="),
    @("80×81=", "53×27="),
    @("61×99=", "61×94="),
    @("67×51=", "58×23="),
    @("55×49=", "23×17="),
    @("97×48=", "50×28="),
    @("55×99=", "93×78="),
    @("63×44=", "44×26="),
    @("85×21=", "38×43="),
    @("63×21=", "71×41="),
    @("88×34=", "67×82="),
    @("51×90=", "48×72="),
    @("23×33=", "47×57="),
    @("45×69=", "55×51="),
    @("40×90=", "41×47="),
    @("37×87=", "73×72="),
    @("32×96=", "72×96="),
    @("64×70=", "85×76="),
    @("21×17=", "58×51="),
    @("75×44=", "65×18="),
    @("55×64=", "47×77=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
